# Merge-driven data wrangle: move the GBR2 rows from Sheet1 to Sheet2,
# and shift the ITA1 / CHN1 / LUX1 rows up to fill the gap on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: append the GBR2 rows (previously Sheet1 rows 13-14) as rows 3-4 ---
$ws2.Range("A3").Value = "GBR2"
$ws2.Range("B3").Value = "region"
$ws2.Range("C3").Value = "data/derived/UK/GBR2_regions.RDS"
$ws2.Range("D3").Value = "marginal"
$ws2.Range("E3").Value = "linelist"

$ws2.Range("A4").Value = "GBR2"
$ws2.Range("B4").Value = "ageband"
$ws2.Range("C4").Value = "data/derived/UK/GBR2_agebands.RDS"
$ws2.Range("D4").Value = "marginal"
$ws2.Range("E4").Value = "linelist"

# --- Sheet1: shift ITA1 / CHN1 / LUX1 rows up into rows 13-16, then clear 17-18 ---
$ws1.Range("A13").Value = "ITA1"
$ws1.Range("B13").Value = "region"
$ws1.Range("C13").Value = "data/derived/ITA/ITA_regions.RDS"
$ws1.Range("D13").Value = "marginal"
$ws1.Range("E13").Value = "aggregate"

$ws1.Range("A14").Value = "ITA1"
$ws1.Range("B14").Value = "ageband"
$ws1.Range("C14").Value = "data/derived/ITA/ITA_agebands.RDS"
$ws1.Range("D14").Value = "marginal"
$ws1.Range("E14").Value = "aggregate"

$ws1.Range("A15").Value = "CHN1"
$ws1.Range("B15").Value = "ageband"
$ws1.Range("C15").Value = "data/derived/CHN/CHN_agebands.RDS"
$ws1.Range("D15").Value = "marginal"
$ws1.Range("E15").Value = "aggregate"

$ws1.Range("A16").Value = "LUX1"
$ws1.Range("B16").Value = "ageband"
$ws1.Range("C16").Value = "data/derived/LUX/LUX_agebands.RDS"
$ws1.Range("D16").Value = "marginal"
$ws1.Range("E16").Value = "aggregate"

$ws1.Range("A17:E18").ClearContents()

# --- selections to mirror the saved view state ---
[void]$ws1.Range("A13:XFD14").Select()
[void]$ws2.Range("A3:XFD4").Select()
[void]$ws1.Activate()
